$wb = $excel.ActiveWorkbook
$ws12 = $wb.Worksheets.Item(13)
$ws12.Range("A10").Value = "MARKER_TEST"
